# Update UI Login for APP
# Adds a new "Icon :" link entry (with bold URL) to Sheet2, row 9,
# and moves the active selection down to B13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# New row 9 / cell B9: rich text "Icon : " (regular) + URL (bold)
$prefix = "Icon : "
$url = "https://www.freepik.com/icon/cancel_18249285#fromView=search&page=1&position=53&uuid=f2f2c95e-2154-4e04-85b1-eb58e78ed2be"

$cell = $ws.Range("B9")
$cell.Value = $prefix + $url

$boldChars = $cell.Characters($prefix.Length + 1, $url.Length)
$boldChars.Font.Bold = $true

# Move the active selection to B13, matching the post-edit workbook state
[void]$ws.Range("B13").Select()
